$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text so values like "0.400" or "0.1000"
# keep their exact digits instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '65.752.36'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.677.76'
$ws.Range('E3').Value = '  -0.99%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '600.82'
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').Value = '156.73'
$ws.Range('E6').Value = '  -0.89%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.618'
$ws.Range('E8').Value = '  +4.97%  '
$ws.Range('E9').Value = '  +3.98%  '
$ws.Range('E10').Value = '  -1.74%  '
$ws.Range('D11').Value = '0.400'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '29.29'
$ws.Range('E13').Value = '  -3.93%  '
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('D15').Value = '3.157.39'
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '65.609.18'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '2.675.50'
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('E18').Value = '  +1.70%  '
$ws.Range('D19').Value = '4.80'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('D21').Value = '352.73'
$ws.Range('E21').Value = '  -1.94%  '
$ws.Range('D23').Value = '69.87'
$ws.Range('E23').Value = '  -1.55%  '
$ws.Range('E24').Value = '  +5.27%  '
$ws.Range('E25').Value = '  -2.55%  '
$ws.Range('D26').Value = '1.65'
$ws.Range('E26').Value = '  -1.18%  '
$ws.Range('D27').Value = '0.168'
$ws.Range('E27').Value = '  -1.98%  '
$ws.Range('E28').Value = '  -6.10%  '
$ws.Range('D29').Value = '8.06'
$ws.Range('E29').Value = '  -4.34%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = '534.02'
$ws.Range('E31').Value = '  -2.00%  '
$ws.Range('E32').Value = '  -3.09%  '
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('E34').Value = '  +1.62%  '
$ws.Range('E35').Value = '  -4.32%  '
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('D37').Value = '20.62'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').Value = '159.54'
$ws.Range('E38').Value = '  -2.28%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  -2.84%  '
$ws.Range('D42').Value = '163.77'
$ws.Range('E42').Value = '  -5.25%  '
$ws.Range('D43').Value = '4.13'
$ws.Range('E43').Value = '  -1.54%  '
$ws.Range('E44').Value = '  +1.83%  '
$ws.Range('E45').Value = '  -1.31%  '
$ws.Range('D46').Value = '22.82'
$ws.Range('E46').Value = '  -3.35%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').Value = '0.640'
$ws.Range('E47').Value = '  -2.53%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').Value = '0.0258'
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('E49').Value = '  +14.92%  '
$ws.Range('D50').Value = '20.27'
$ws.Range('E50').Value = '  -3.76%  '
$ws.Range('D51').Value = '0.1000'
$ws.Range('E51').Value = '  +0.72%  '

# Restore default (no explicit) style on the Price column, matching the source formatting.
$ws.Range("D2:D51").Style = "Normal"
